$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (Sheet1 -> Andre Russell)
$ws.Name = "Andre Russell"

# Insert a new first column for "matchNo" (shifts teamName..result right by one)
$ws.Range("A1").EntireColumn.Insert()

# Header row
$ws.Cells.Item(1,1).Value = 'matchNo'
$ws.Cells.Item(1,2).Value = 'teamName'
$ws.Cells.Item(1,3).Value = 'batterName'
$ws.Cells.Item(1,4).Value = 'states'
$ws.Cells.Item(1,5).Value = 'runs'
$ws.Cells.Item(1,6).Value = 'balls'
$ws.Cells.Item(1,7).Value = 'fours'
$ws.Cells.Item(1,8).Value = 'sixes'
$ws.Cells.Item(1,9).Value = 'sr'
$ws.Cells.Item(1,10).Value = 'opponentTeamName'
$ws.Cells.Item(1,11).Value = 'venue'
$ws.Cells.Item(1,12).Value = 'date'
$ws.Cells.Item(1,13).Value = 'result'

# Data rows
# Row 2
$ws.Cells.Item(2,1).Value = '10th'
$ws.Cells.Item(2,2).Value = 'Kolkata Knight Riders'
$ws.Cells.Item(2,3).Value = 'Andre Russell'
$ws.Cells.Item(2,4).Value = 'b Patel'
$ws.Cells.Item(2,5).NumberFormat = "@"
$ws.Cells.Item(2,5).Value = '31'
$ws.Cells.Item(2,6).NumberFormat = "@"
$ws.Cells.Item(2,6).Value = '20'
$ws.Cells.Item(2,7).NumberFormat = "@"
$ws.Cells.Item(2,7).Value = '3'
$ws.Cells.Item(2,8).NumberFormat = "@"
$ws.Cells.Item(2,8).Value = '2'
$ws.Cells.Item(2,9).NumberFormat = "@"
$ws.Cells.Item(2,9).Value = '155.00'
$ws.Cells.Item(2,10).Value = 'Royal Challengers Bangalore'
$ws.Cells.Item(2,11).Value = 'Chennai'
$ws.Cells.Item(2,12).Value = 'April 18'
$ws.Cells.Item(2,13).Value = 'RCB won by 38 runs'
# Row 3
$ws.Cells.Item(3,1).Value = '31st'
$ws.Cells.Item(3,2).Value = 'Kolkata Knight Riders'
$ws.Cells.Item(3,3).Value = 'Andre Russell'
$ws.Cells.Item(3,4).Formula = '=""'
$ws.Cells.Item(3,5).NumberFormat = "@"
$ws.Cells.Item(3,5).Value = '0'
$ws.Cells.Item(3,6).NumberFormat = "@"
$ws.Cells.Item(3,6).Value = '0'
$ws.Cells.Item(3,7).NumberFormat = "@"
$ws.Cells.Item(3,7).Value = '0'
$ws.Cells.Item(3,8).NumberFormat = "@"
$ws.Cells.Item(3,8).Value = '0'
$ws.Cells.Item(3,9).NumberFormat = "@"
$ws.Cells.Item(3,9).Value = '-'
$ws.Cells.Item(3,10).Value = 'Royal Challengers Bangalore'
$ws.Cells.Item(3,11).Value = 'Abu Dhabi'
$ws.Cells.Item(3,12).Value = 'September 20'
$ws.Cells.Item(3,13).Value = 'KKR won by 9 wickets (with 60 balls remaining)'
# Row 4
$ws.Cells.Item(4,1).Value = '25th'
$ws.Cells.Item(4,2).Value = 'Kolkata Knight Riders'
$ws.Cells.Item(4,3).Value = 'Andre Russell'
$ws.Cells.Item(4,4).Formula = '=""'
$ws.Cells.Item(4,5).NumberFormat = "@"
$ws.Cells.Item(4,5).Value = '45'
$ws.Cells.Item(4,6).NumberFormat = "@"
$ws.Cells.Item(4,6).Value = '27'
$ws.Cells.Item(4,7).NumberFormat = "@"
$ws.Cells.Item(4,7).Value = '2'
$ws.Cells.Item(4,8).NumberFormat = "@"
$ws.Cells.Item(4,8).Value = '4'
$ws.Cells.Item(4,9).NumberFormat = "@"
$ws.Cells.Item(4,9).Value = '166.66'
$ws.Cells.Item(4,10).Value = 'Delhi Capitals'
$ws.Cells.Item(4,11).Value = 'Ahmedabad'
$ws.Cells.Item(4,12).Value = 'April 29'
$ws.Cells.Item(4,13).Value = 'Capitals won by 7 wickets (with 21 balls remaining)'
# Row 5
$ws.Cells.Item(5,1).Value = '15th'
$ws.Cells.Item(5,2).Value = 'Kolkata Knight Riders'
$ws.Cells.Item(5,3).Value = 'Andre Russell'
$ws.Cells.Item(5,4).Value = 'b Curran'
$ws.Cells.Item(5,5).NumberFormat = "@"
$ws.Cells.Item(5,5).Value = '54'
$ws.Cells.Item(5,6).NumberFormat = "@"
$ws.Cells.Item(5,6).Value = '22'
$ws.Cells.Item(5,7).NumberFormat = "@"
$ws.Cells.Item(5,7).Value = '3'
$ws.Cells.Item(5,8).NumberFormat = "@"
$ws.Cells.Item(5,8).Value = '6'
$ws.Cells.Item(5,9).NumberFormat = "@"
$ws.Cells.Item(5,9).Value = '245.45'
$ws.Cells.Item(5,10).Value = 'Chennai Super Kings'
$ws.Cells.Item(5,11).Value = 'Wankhede'
$ws.Cells.Item(5,12).Value = 'April 21'
$ws.Cells.Item(5,13).Value = 'Super Kings won by 18 runs'
# Row 6
$ws.Cells.Item(6,1).Value = '18th'
$ws.Cells.Item(6,2).Value = 'Kolkata Knight Riders'
$ws.Cells.Item(6,3).Value = 'Andre Russell'
$ws.Cells.Item(6,4).Value = 'c Miller b Morris'
$ws.Cells.Item(6,5).NumberFormat = "@"
$ws.Cells.Item(6,5).Value = '9'
$ws.Cells.Item(6,6).NumberFormat = "@"
$ws.Cells.Item(6,6).Value = '7'
$ws.Cells.Item(6,7).NumberFormat = "@"
$ws.Cells.Item(6,7).Value = '0'
$ws.Cells.Item(6,8).NumberFormat = "@"
$ws.Cells.Item(6,8).Value = '1'
$ws.Cells.Item(6,9).NumberFormat = "@"
$ws.Cells.Item(6,9).Value = '128.57'
$ws.Cells.Item(6,10).Value = 'Rajasthan Royals'
$ws.Cells.Item(6,11).Value = 'Wankhede'
$ws.Cells.Item(6,12).Value = 'April 24'
$ws.Cells.Item(6,13).Value = 'Royals won by 6 wickets (with 7 balls remaining)'
# Row 7
$ws.Cells.Item(7,1).Value = '21st'
$ws.Cells.Item(7,2).Value = 'Kolkata Knight Riders'
$ws.Cells.Item(7,3).Value = 'Andre Russell'
$ws.Cells.Item(7,4).Value = 'run out (Arshdeep Singh/†Rahul)'
$ws.Cells.Item(7,5).NumberFormat = "@"
$ws.Cells.Item(7,5).Value = '10'
$ws.Cells.Item(7,6).NumberFormat = "@"
$ws.Cells.Item(7,6).Value = '9'
$ws.Cells.Item(7,7).NumberFormat = "@"
$ws.Cells.Item(7,7).Value = '2'
$ws.Cells.Item(7,8).NumberFormat = "@"
$ws.Cells.Item(7,8).Value = '0'
$ws.Cells.Item(7,9).NumberFormat = "@"
$ws.Cells.Item(7,9).Value = '111.11'
$ws.Cells.Item(7,10).Value = 'Punjab Kings'
$ws.Cells.Item(7,11).Value = 'Ahmedabad'
$ws.Cells.Item(7,12).Value = 'April 26'
$ws.Cells.Item(7,13).Value = 'KKR won by 5 wickets (with 20 balls remaining)'
# Row 8
$ws.Cells.Item(8,1).Value = '3rd'
$ws.Cells.Item(8,2).Value = 'Kolkata Knight Riders'
$ws.Cells.Item(8,3).Value = 'Andre Russell'
$ws.Cells.Item(8,4).Value = 'c Pandey b Rashid Khan'
$ws.Cells.Item(8,5).NumberFormat = "@"
$ws.Cells.Item(8,5).Value = '5'
$ws.Cells.Item(8,6).NumberFormat = "@"
$ws.Cells.Item(8,6).Value = '5'
$ws.Cells.Item(8,7).NumberFormat = "@"
$ws.Cells.Item(8,7).Value = '1'
$ws.Cells.Item(8,8).NumberFormat = "@"
$ws.Cells.Item(8,8).Value = '0'
$ws.Cells.Item(8,9).NumberFormat = "@"
$ws.Cells.Item(8,9).Value = '100.00'
$ws.Cells.Item(8,10).Value = 'Sunrisers Hyderabad'
$ws.Cells.Item(8,11).Value = 'Chennai'
$ws.Cells.Item(8,12).Value = 'April 11'
$ws.Cells.Item(8,13).Value = 'KKR won by 10 runs'
# Row 9
$ws.Cells.Item(9,1).Value = '5th'
$ws.Cells.Item(9,2).Value = 'Kolkata Knight Riders'
$ws.Cells.Item(9,3).Value = 'Andre Russell'
$ws.Cells.Item(9,4).Value = 'c & b Boult'
$ws.Cells.Item(9,5).NumberFormat = "@"
$ws.Cells.Item(9,5).Value = '9'
$ws.Cells.Item(9,6).NumberFormat = "@"
$ws.Cells.Item(9,6).Value = '15'
$ws.Cells.Item(9,7).NumberFormat = "@"
$ws.Cells.Item(9,7).Value = '1'
$ws.Cells.Item(9,8).NumberFormat = "@"
$ws.Cells.Item(9,8).Value = '0'
$ws.Cells.Item(9,9).NumberFormat = "@"
$ws.Cells.Item(9,9).Value = '60.00'
$ws.Cells.Item(9,10).Value = 'Mumbai Indians'
$ws.Cells.Item(9,11).Value = 'Chennai'
$ws.Cells.Item(9,12).Value = 'April 13'
$ws.Cells.Item(9,13).Value = 'Mumbai won by 10 runs'
# Row 10
$ws.Cells.Item(10,1).Value = '38th'
$ws.Cells.Item(10,2).Value = 'Kolkata Knight Riders'
$ws.Cells.Item(10,3).Value = 'Andre Russell'
$ws.Cells.Item(10,4).Value = 'b Thakur'
$ws.Cells.Item(10,5).NumberFormat = "@"
$ws.Cells.Item(10,5).Value = '20'
$ws.Cells.Item(10,6).NumberFormat = "@"
$ws.Cells.Item(10,6).Value = '15'
$ws.Cells.Item(10,7).NumberFormat = "@"
$ws.Cells.Item(10,7).Value = '2'
$ws.Cells.Item(10,8).NumberFormat = "@"
$ws.Cells.Item(10,8).Value = '1'
$ws.Cells.Item(10,9).NumberFormat = "@"
$ws.Cells.Item(10,9).Value = '133.33'
$ws.Cells.Item(10,10).Value = 'Chennai Super Kings'
$ws.Cells.Item(10,11).Value = 'Abu Dhabi'
$ws.Cells.Item(10,12).Value = 'September 26'
$ws.Cells.Item(10,13).Value = 'Super Kings won by 2 wickets'
